$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "56.386.04"
Set-TextValue $ws.Range("E2") "  -0.51%  "
Set-TextValue $ws.Range("D3") "2.312.90"
Set-TextValue $ws.Range("E3") "  -0.46%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "511.28"
Set-TextValue $ws.Range("E5") "  -1.81%  "
Set-TextValue $ws.Range("D6") "130.90"
Set-TextValue $ws.Range("E6") "  -2.97%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.52%  "
Set-TextValue $ws.Range("D8") "0.532"
Set-TextValue $ws.Range("E8") "  -1.02%  "
Set-TextValue $ws.Range("E9") "  -3.72%  "
Set-TextValue $ws.Range("E10") "  -0.15%  "
Set-TextValue $ws.Range("D11") "5.24"
Set-TextValue $ws.Range("E11") "  -0.19%  "
Set-TextValue $ws.Range("D12") "0.335"
Set-TextValue $ws.Range("E12") "  -2.12%  "
Set-TextValue $ws.Range("D13") "2.733.58"
Set-TextValue $ws.Range("E13") "  -0.24%  "
Set-TextValue $ws.Range("D14") "23.42"
Set-TextValue $ws.Range("E14") "  -1.95%  "
Set-TextValue $ws.Range("D15") "56.385.53"
Set-TextValue $ws.Range("E15") "  -0.62%  "
Set-TextValue $ws.Range("E16") "  -2.15%  "
Set-TextValue $ws.Range("D17") "2.320.21"
Set-TextValue $ws.Range("E17") "  -0.29%  "
Set-TextValue $ws.Range("D18") "10.35"
Set-TextValue $ws.Range("E18") "  -1.15%  "
Set-TextValue $ws.Range("D19") "326.25"
Set-TextValue $ws.Range("E19") "  +0.93%  "
Set-TextValue $ws.Range("D20") "4.12"
Set-TextValue $ws.Range("E20") "  -2.46%  "
Set-TextValue $ws.Range("E21") "  +1.93%  "
Set-TextValue $ws.Range("D22") "0.998"
Set-TextValue $ws.Range("E22") "  -0.14%  "
Set-TextValue $ws.Range("D23") "61.24"
Set-TextValue $ws.Range("E23") "  +1.11%  "
Set-TextValue $ws.Range("D24") "8.55"
Set-TextValue $ws.Range("E24") "  +7.97%  "
Set-TextValue $ws.Range("D25") "0.162"
Set-TextValue $ws.Range("E25") "  -1.00%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.49%  "
Set-TextValue $ws.Range("E27") "  +0.05%  "
Set-TextValue $ws.Range("D28") "167.57"
Set-TextValue $ws.Range("E28") "  -1.33%  "
Set-TextValue $ws.Range("E29") "  -2.86%  "
Set-TextValue $ws.Range("D30") "0.0₃0715"
Set-TextValue $ws.Range("E30") "  -4.34%  "
Set-TextValue $ws.Range("D31") "6.08"
Set-TextValue $ws.Range("E31") "  -1.63%  "
Set-TextValue $ws.Range("D32") "18.28"
Set-TextValue $ws.Range("E32") "  -0.41%  "
Set-TextValue $ws.Range("E33") "  -0.02%  "
Set-TextValue $ws.Range("E34") "  +0.63%  "
Set-TextValue $ws.Range("E35") "  +0.15%  "
Set-TextValue $ws.Range("D36") "3.91"
Set-TextValue $ws.Range("E36") "  -3.52%  "
Set-TextValue $ws.Range("D37") "0.877"
Set-TextValue $ws.Range("E37") "  -5.15%  "
Set-TextValue $ws.Range("D38") "38.50"
Set-TextValue $ws.Range("E38") "  +1.44%  "
Set-TextValue $ws.Range("D39") "1.55"
Set-TextValue $ws.Range("E39") "  -0.21%  "
Set-TextValue $ws.Range("D40") "148.73"
Set-TextValue $ws.Range("E40") "  +7.82%  "
Set-TextValue $ws.Range("D41") "0.372"
Set-TextValue $ws.Range("E41") "  -1.88%  "
Set-TextValue $ws.Range("E42") "  -0.74%  "
Set-TextValue $ws.Range("D43") "274.24"
Set-TextValue $ws.Range("E43") "  -0.81%  "
Set-TextValue $ws.Range("D44") "4.97"
Set-TextValue $ws.Range("E44") "  -5.49%  "
Set-TextValue $ws.Range("D45") "0.0925"
Set-TextValue $ws.Range("E45") "  -0.84%  "
Set-TextValue $ws.Range("D46") "0.0493"
Set-TextValue $ws.Range("E46") "  -2.52%  "
Set-TextValue $ws.Range("D47") "0.553"
Set-TextValue $ws.Range("E47") "  -1.79%  "
Set-TextValue $ws.Range("D48") "18.08"
Set-TextValue $ws.Range("E48") "  +0.71%  "
Set-TextValue $ws.Range("E49") "  -0.08%  "
Set-TextValue $ws.Range("E50") "  -1.82%  "
Set-TextValue $ws.Range("D51") "16.94"
Set-TextValue $ws.Range("E51") "  +0.54%  "
